$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Individual cell value updates per diff
$ws.Range("J2").Value = 2.62
$ws.Range("K2").Value = 2.37
$ws.Range("Y2").Value = 1.5
$ws.Range("Z2").Value = 2.37
$ws.Range("AH2").Value = 8
$ws.Range("AN2").Value = 13
$ws.Range("AP2").Value = 26
$ws.Range("O3").Value = 1.29
$ws.Range("P3").Value = 3.75
$ws.Range("Y3").Value = 1.72
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 2.5
$ws.Range("L4").Value = 3.2
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10
$ws.Range("AD4").Value = 29
$ws.Range("AO4").Value = 29
$ws.Range("AP4").Value = 23
$ws.Range("J6").Value = 2.67
$ws.Range("L6").Value = 3.7
$ws.Range("O6").Value = 1.26
$ws.Range("P6").Value = 3.5
$ws.Range("U6").Value = 2.82
$ws.Range("AA6").Value = 8.5
$ws.Range("AL6").Value = 10.75
$ws.Range("AM6").Value = 17.5
$ws.Range("Q11").Value = 2
$ws.Range("R11").Value = 1.85
$ws.Range("U11").Value = 3.4
$ws.Range("V11").Value = 1.3
$ws.Range("G12").Value = 2.8
$ws.Range("I12").Value = 2.55
$ws.Range("J12").Value = 3.75
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 7
$ws.Range("U12").Value = 5
$ws.Range("V12").Value = 1.17
$ws.Range("Y12").Value = 2.05
$ws.Range("Z12").Value = 1.7
$ws.Range("AG12").Value = 7
$ws.Range("AI12").Value = 19
$ws.Range("AK12").Value = 501
$ws.Range("AL12").Value = 6.5
$ws.Range("AR12").Value = 1.85
$ws.Range("AS12").Value = 2
$ws.Range("G14").Value = 1.39
$ws.Range("I15").Value = 1.69
$ws.Range("G16").Value = 1.1
$ws.Range("I17").Value = 1.39
$ws.Range("AK17").Value = 450
$ws.Range("G19").Value = 1.13
$ws.Range("Q19").Value = 1.3
$ws.Range("G20").Value = 3.4
$ws.Range("I20").Value = 2.1
$ws.Range("J20").Value = 3.75
$ws.Range("L20").Value = 2.63
$ws.Range("M20").Value = 1.05
$ws.Range("N20").Value = 8.5
$ws.Range("AA20").Value = 12
$ws.Range("AB20").Value = 19
$ws.Range("AC20").Value = 13
$ws.Range("AE20").Value = 26
$ws.Range("AF20").Value = 34
$ws.Range("AL20").Value = 8.5
$ws.Range("AN20").Value = 9
$ws.Range("AO20").Value = 19
$ws.Range("G22").Value = 1.45
$ws.Range("U23").Value = 2.62

# Row 13 (Bengaluru FC vs Chennaiyin) - fill in previously empty odds columns G:AQ
$ws.Range("G13").Value = 1.91
$ws.Range("H13").Value = 3.5
$ws.Range("I13").Value = 3.6
$ws.Range("J13").Value = 2.6
$ws.Range("K13").Value = 2.25
$ws.Range("L13").Value = 4
$ws.Range("M13").Value = 1.04
$ws.Range("N13").Value = 13
$ws.Range("O13").Value = 1.22
$ws.Range("P13").Value = 4
$ws.Range("Q13").Value = 1.73
$ws.Range("R13").Value = 2.08
$ws.Range("S13").Value = 2
$ws.Range("T13").Value = 1.8
$ws.Range("U13").Value = 2.75
$ws.Range("V13").Value = 1.4
$ws.Range("W13").Value = 1.36
$ws.Range("X13").Value = 3
$ws.Range("Y13").Value = 1.67
$ws.Range("Z13").Value = 2.1
$ws.Range("AA13").Value = 8.5
$ws.Range("AB13").Value = 10
$ws.Range("AC13").Value = 8.5
$ws.Range("AD13").Value = 17
$ws.Range("AE13").Value = 15
$ws.Range("AF13").Value = 23
$ws.Range("AG13").Value = 12
$ws.Range("AH13").Value = 7
$ws.Range("AI13").Value = 13
$ws.Range("AJ13").Value = 41
$ws.Range("AK13").Value = 151
$ws.Range("AL13").Value = 13
$ws.Range("AM13").Value = 21
$ws.Range("AN13").Value = 13
$ws.Range("AO13").Value = 41
$ws.Range("AP13").Value = 29
$ws.Range("AQ13").Value = 34
